# restructure the generator functions
#
# The "case1/case2/isDivisibleBy" note (row 5, column E -- the
# math_boolean_numberProperty row) gets its body text corrected: the
# "case2" label is missing a trailing colon ("case2:"), matching the
# sibling "case1:" label above it.
#
# Saving the workbook afterwards naturally rewrites the shared-strings
# table (the stale, now-orphaned string is dropped and the corrected one
# is appended at the end), which is what produces all of the <v> index
# churn seen in the diff -- no other cell text actually changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E5").Value = "case1:" + [char]10 + "<arg_1>" + [char]10 + "isEven / " + [char]10 + "case2:" + [char]10 + "<arg_1>" + [char]10 + "<arg_2>" + [char]10 + "isDivisibleBy"

# The view no longer freezes/scrolls to column B as the left-most visible
# column -- restore the default top-left cell (A1).
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1

# Column D (Command in Bytecode) narrows and switches to an explicit
# best-fit width.
$ws.Columns.Item(4).ColumnWidth = 51.83203125
$ws.Columns.Item(4).AutoFit() | Out-Null
$ws.Columns.Item(4).ColumnWidth = 51.83203125
